# Re-purpose room C001 (row 2) from a "Recreation" room of unknown
# ("nil") capacity into a "large classroom" seating 120, and bump the
# existing C002 room's (row 3) capacity from 116 to 120 as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 used to describe room C001 as a "Recreation" room with capacity "nil"
# and facility "None". It is retyped to describe it as a "large classroom"
# with a capacity of 120 (facility stays "None").
$ws.Range("B2").Value = "large classroom"
$ws.Range("C2").Value = 120

# Row 3 (room C002, "large classroom") had capacity 116; update it to 120.
$ws.Range("C3").Value = 120

# Reflect the selected cell as left in the saved workbook.
$ws.Range("C1").Select()
